$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999998762278808,
    0.99417385683522586,
    0.97318704609608653,
    0.96357332751712432,
    0.95437536625733133,
    0.93387042941657294,
    0.93130237090660917,
    0.92702768374521705,
    0.92420760498970633,
    0.92231535823632316,
    0.92203677776028337,
    0.92176289694234448,
    0.92349301805840511,
    0.92548956753444866,
    0.92841134910993861,
    0.93229841170143279,
    0.94340859884903017,
    0.94229973203462269,
    0.99117740329712856,
    0.97636358109617694,
    0.9740244442037177,
    0.97275993788555026,
    0.96615157126888807,
    0.95312997845951997,
    0.94667288152647799,
    0.93128198511900484,
    0.92792165869992704,
    0.91517674648943581,
    0.90636574452695196,
    0.90233167657095925,
    0.90335196343741386,
    0.90489170464442381,
    0.90437171199059951
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
